# Apply "Added Test Data for Partial Circuit Cabling" changes
$wb = $excel.ActiveWorkbook

$wsDeviceInfo = $wb.Worksheets.Item("DeviceInfo")
$wsDSP = $wb.Worksheets.Item("DSPChannelMap")

# Update DSPChannelMap data values (D5:D13 and B11:B13)
$wsDSP.Range("D5").Value = 10
$wsDSP.Range("D6").Value = 11
$wsDSP.Range("D7").Value = 12
$wsDSP.Range("D8").Value = 13
$wsDSP.Range("D9").Value = 14
$wsDSP.Range("D10").Value = 15
$wsDSP.Range("B11").Value = 16
$wsDSP.Range("B12").Value = 17
$wsDSP.Range("B13").Value = 18

# Update selections to match final state
$wsDeviceInfo.Range("I14").Select()
$wsDSP.Range("B13").Select()

# Make DSPChannelMap the active sheet (so it becomes tabSelected and activeTab)
$wsDSP.Activate()
